# "Add FO, Budgets and Winner FK"
# Adds a new timeline entry (row 30) describing the SaveAndLoad / date-time
# formatting / ForeignKey-algorithm / multi-field-incl.-FK work, and moves
# the sheet's active selection on to B31 (the next empty row) as the author
# left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newNote = "Работа по созданию функционала загрузки данных (Функция SaveAndLoad, Формат даты и времени, Изменения алгоритма для ForeignKey, Настройка нескольких полей включая FK)"

# Row 30: A = description, B = hours spent, C = date completed,
# D already carries the shared formula "=B*$B$1" and recalculates itself.
$ws.Range("A30").Value = $newNote
$ws.Range("B30").Value = 3

# Give C30 the same date number-format as the rest of the "Дата выполнения"
# column (copy format from C29) before writing the date value itself.
$ws.Range("C29").Copy()
$ws.Range("C30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C30").Value = (Get-Date -Year 2019 -Month 4 -Day 30).Date

# Move the cursor/selection on to the next row, matching the author's
# final on-screen state.
$ws.Activate()
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 2
